# Apply the updated cryptocurrency price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.061.36"
$ws.Range("E2").Value = "  +2.83%  "
$ws.Range("D3").Value = "'2.654.45"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'595.69"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").Value = "'156.11"
$ws.Range("E6").Value = "  +3.84%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("E9").Value = "  +7.76%  "
$ws.Range("E10").Value = "  +4.19%  "
$ws.Range("D11").Value = "'5.79"
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("D13").Value = "'29.00"
$ws.Range("E13").Value = "  +5.12%  "
$ws.Range("E14").Value = "  +20.44%  "
$ws.Range("D15").Value = "'3.131.26"
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("D16").Value = "'64.987.35"
$ws.Range("E16").Value = "  +3.05%  "
$ws.Range("D17").Value = "'2.651.45"
$ws.Range("E17").Value = "  +3.46%  "
$ws.Range("D18").Value = "'12.62"
$ws.Range("E18").Value = "  +3.73%  "
$ws.Range("E19").Value = "  +1.91%  "
$ws.Range("D20").Value = "'354.20"
$ws.Range("E20").Value = "  +2.98%  "
$ws.Range("E21").Value = "  +6.56%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "'68.19"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("D25").Value = "'9.51"
$ws.Range("E25").Value = "  +3.77%  "
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").Value = "'8.17"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0₃0946"
$ws.Range("E30").Value = "  +10.75%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.12"
$ws.Range("E31").Value = "  +5.11%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'522.57"
$ws.Range("E32").Value = "  -6.26%  "
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("E34").Value = "  +8.45%  "
$ws.Range("D35").Value = "'6.37"
$ws.Range("E35").Value = "  +4.80%  "
$ws.Range("D36").Value = "'0.427"
$ws.Range("E36").Value = "  +3.90%  "
$ws.Range("D37").Value = "'165.19"
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("D38").Value = "'20.27"
$ws.Range("E38").Value = "  +3.90%  "
$ws.Range("E39").Value = "  +5.49%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "'42.20"
$ws.Range("E42").Value = "  +6.52%  "
$ws.Range("D43").Value = "'165.51"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("E44").Value = "  +3.24%  "
$ws.Range("D45").Value = "'0.0618"
$ws.Range("E45").Value = "  +6.25%  "
$ws.Range("D46").Value = "'22.99"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("E47").Value = "  +4.58%  "
$ws.Range("E48").Value = "  +3.69%  "
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("D51").Value = "'19.51"
$ws.Range("E51").Value = "  +1.81%  "
